# Updated cryptos list on Sat Aug 26 03:55:31 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (col D) and "Volume(1h)" (col E) columns for each
# coin row (rows 2-51) with the latest scraped figures. Price strings that
# look like plain numbers are written with a leading apostrophe so Excel
# keeps them as text (matching the original inlineStr cells) instead of
# silently re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PV($Row, $Price, $Volume) {
    if ($Price) {
        $text = $Price
        if ($text -match '^[0-9]+(\.[0-9]+)?$') {
            $text = "'" + $text
        }
        $ws.Cells.Item($Row, 4).Value = $text
    }
    $ws.Cells.Item($Row, 5).Value = "  $Volume  "
}

Set-PV 2  "26.144.86" "-0.31%"
Set-PV 3  "1.657.20"  "-0.32%"
Set-PV 4  $null       "-0.33%"
Set-PV 5  "218.17"    "+0.14%"
Set-PV 6  "0.5287"    "+1.15%"
Set-PV 7  $null       "-0.27%"
Set-PV 8  "0.2611"    "-1.23%"
Set-PV 9  "0.06350"   "+1.08%"
Set-PV 10 "20.47"     "-1.57%"
Set-PV 11 "0.07782"   "+0.64%"
Set-PV 12 $null       "+1.85%"
Set-PV 13 "1.661.25"  "-0.08%"
Set-PV 14 "0.5502"    "+1.26%"
Set-PV 15 "0.0₅8224"  "+1.02%"
Set-PV 16 "65.49"     "+1.52%"
Set-PV 17 "26.153.09" "-0.37%"
Set-PV 18 $null       "-0.31%"
Set-PV 19 "4.584"     "-1.27%"
Set-PV 20 "192.57"    "-0.23%"
Set-PV 21 $null       "+0.44%"
Set-PV 22 $null       "+0.00%"
Set-PV 23 "1.004"     "-0.47%"
Set-PV 24 $null       "+1.40%"
Set-PV 25 "0.1248"    "+1.55%"
Set-PV 26 "7.286"     "+1.78%"
Set-PV 27 $null       "+0.68%"
Set-PV 28 "1.441"     "+2.02%"
Set-PV 29 "0.05936"   "-2.39%"
Set-PV 30 "1.279"     "+0.11%"
Set-PV 31 "3.529"     "-1.29%"
Set-PV 32 "3.269"     "+0.35%"
Set-PV 33 "1.587"     "-2.07%"
Set-PV 34 "0.9560"    "-0.89%"
Set-PV 35 "2.793"     "+0.35%"
Set-PV 36 "2.410"     "-0.72%"
Set-PV 37 $null       "+0.57%"
Set-PV 38 $null       "+1.90%"
Set-PV 39 "5.806"     "-2.89%"
Set-PV 40 "0.8483"    "-0.74%"
Set-PV 41 $null       "-0.23%"
Set-PV 42 "103.28"    "+3.03%"
Set-PV 43 "1.026.55"  "+1.37%"
Set-PV 44 "1.801.97"  "-0.15%"
Set-PV 45 "57.44"     "+0.65%"
Set-PV 46 "1.007"     "-0.13%"
Set-PV 47 "1.486"     "+0.75%"
Set-PV 48 "0.4294"    "+1.81%"
Set-PV 49 $null       "-0.45%"
Set-PV 50 "7.819"     "-1.89%"
Set-PV 51 "0.09726"   "+0.27%"
